$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing measurement for row 16 (I16), which also
# recalculates K16 (=I16/J16) from 0 to its real ratio.
$ws.Range("I16").Value = 488

# Bring over the formatting of row 16 onto the new row 17 (styles only),
# mirroring a copy/paste of the row before editing its contents.
$ws.Range("B16:L16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# New data row appended below the existing table.
$ws.Range("B17").Value = 42687
$ws.Range("C17").Value = 900
$ws.Range("D17").Value = 2.1
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 233
$ws.Range("J17").Formula = "=F17*E17*G17/H17"
$ws.Range("K17").Formula = "=I17/J17"
$ws.Range("L17").Formula = "=G17*F17*E17"

# Leave the sheet scrolled/selected the same way the author ended up with.
$ws.Range("I18").Select() | Out-Null
